# Insert a new price record at row 291 for "Ajo / Chino / Primera" (Vega
# Central Mapocho de Santiago). All the existing rows from 291 downward
# shift down by one (to 292..360); the new row carries an updated date
# and price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 291..359 down to 292..360, leaving a blank row 291.
$ws.Rows.Item(291).Insert()

# Populate the newly inserted row 291 with the new record.
$ws.Cells.Item(291, 1).Value  = 9
$ws.Cells.Item(291, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(291, 3).Value  = "Metropolitana"
$ws.Cells.Item(291, 4).Value  = 45204
$ws.Cells.Item(291, 5).Value  = 13
$ws.Cells.Item(291, 6).Value  = 100112003
$ws.Cells.Item(291, 7).Value  = "Ajo"
$ws.Cells.Item(291, 8).Value  = "Chino"
$ws.Cells.Item(291, 9).Value  = "Primera"
$ws.Cells.Item(291, 10).Value = 340
$ws.Cells.Item(291, 11).Value = 19000
$ws.Cells.Item(291, 12).Value = 20000
$ws.Cells.Item(291, 13).Value = 19500
$ws.Cells.Item(291, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(291, 15).Value = "China"
$ws.Cells.Item(291, 16).Value = 1950
$ws.Cells.Item(291, 17).Value = 10
$ws.Cells.Item(291, 18).Value = "Hortaliza"
